$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 3088
$ws1.Range("F6").Value = 2063
$ws1.Range("F9").Value = 1157
$ws1.Range("F11").Value = 921
$ws1.Range("F12").Value = 78

# Sheet "全部类型"
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F5").Value = 3088
$ws2.Range("F6").Value = 2063
$ws2.Range("F10").Value = 1157
$ws2.Range("F12").Value = 921
$ws2.Range("F13").Value = 78
